# Update the re_rank (G) scores, and fix a couple of mis-ordered rows
# (prolificid/name/gender got shifted by one position within their
# race group) so the ranking table / infobox reads correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Asian group: rows 2-13 ---
$ws.Range("G2").Value2 = 11.35146450363736
$ws.Range("G3").Value2 = 10.46509312749224

# Row 4/5 swap (prolificid + name)
$ws.Range("C4").Value2 = 2
$ws.Range("D4").Value2 = "5f2c1a97a6809c060fec8820"
$ws.Range("E4").Value2 = "Maggie"
$ws.Range("G4").Value2 = 8.345780257992518

$ws.Range("C5").Value2 = 10
$ws.Range("D5").Value2 = "60a71d27a66fac796ad4de6f"
$ws.Range("E5").Value2 = "Jennifer"
$ws.Range("G5").Value2 = 8.335815743434177

$ws.Range("G6").Value2 = 7.149319585641714
$ws.Range("G7").Value2 = 6.428054772178267
$ws.Range("G8").Value2 = 6.387792560901622
$ws.Range("G9").Value2 = 5.358859522459975
$ws.Range("G10").Value2 = 5.046880051605173
$ws.Range("G11").Value2 = 4.123973324417245
$ws.Range("G12").Value2 = 2.400210451344222
$ws.Range("G13").Value2 = 0.1240641252646651

# --- Hispanic group: rows 14-25 ---
$ws.Range("G14").Value2 = 8.43712372573634
$ws.Range("G15").Value2 = 8.085171174491482

# Row 16/17 swap (prolificid + name + gender)
$ws.Range("C16").Value2 = 3
$ws.Range("D16").Value2 = "60ba8ba51a5e0a105396888a"
$ws.Range("E16").Value2 = "Alfredo"
$ws.Range("F16").Value2 = "male"
$ws.Range("G16").Value2 = 7.483015296297952

$ws.Range("C17").Value2 = 8
$ws.Range("D17").Value2 = "5f0142aa1eb1e528e7abce50"
$ws.Range("E17").Value2 = "Valeria"
$ws.Range("F17").Value2 = "female"
$ws.Range("G17").Value2 = 7.18543091573438

$ws.Range("G18").Value2 = 6.391489871176138
$ws.Range("G19").Value2 = 6.387011644639443
$ws.Range("G20").Value2 = 5.012348235563821

# Row 21/22 swap (prolificid + name)
$ws.Range("C21").Value2 = 4
$ws.Range("D21").Value2 = "5e706891c396cc64388ef760"
$ws.Range("E21").Value2 = "Maria"
$ws.Range("G21").Value2 = 3.142757349846526

$ws.Range("C22").Value2 = 1
$ws.Range("D22").Value2 = "5e0adc8f4cac6834756db412"
$ws.Range("E22").Value2 = "Mary"
$ws.Range("G22").Value2 = 3.037402404511541

$ws.Range("G23").Value2 = 2.344227293246886
$ws.Range("G24").Value2 = 1.337899327810482
$ws.Range("G25").Value2 = 0.002543171126171584
